# Generate Report for Handback
#
# The handback transform for the 845cec40-... file failed, so the report
# needs to reflect the failure on every sheet that tracks it (Overview,
# zh-cn, de-de) and record the error detail message on the per-locale
# sheets. The "Error Detail" column (P) also needs to be widened so the
# message is readable.

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

$zhErrorDetail = "Handback file name: u1x3lneo.12n is different with handoff file name: 845cec40-b0a2-4467-9860-bedb8d4912d5.9439c1348d0544edb68af8fc6231f899c48aefe5.zh-cn."
$deErrorDetail = "Handback file name: u1x3lneo.12n is different with handoff file name: 845cec40-b0a2-4467-9860-bedb8d4912d5.9439c1348d0544edb68af8fc6231f899c48aefe5.de-de."

# --- Overview sheet: update the per-locale status for the 845cec40 row ---
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E3").Value = $statusFailed
$overview.Range("F3").Value = $statusFailed

# --- zh-cn sheet: status + error detail for the 845cec40 row, widen col P ---
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusFailed
$zhcn.Range("P3").Value = $zhErrorDetail
$zhcn.Columns.Item(16).ColumnWidth = $zhcn.Columns.Item(1).ColumnWidth

# --- de-de sheet: status + error detail for the 845cec40 row, widen col P ---
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C3").Value = $statusFailed
$dede.Range("P3").Value = $deErrorDetail
$dede.Columns.Item(16).ColumnWidth = $dede.Columns.Item(1).ColumnWidth
